$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.774.00"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.044.63"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "227.58"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "60.19"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D10").Value = "0.0834"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "2.346.31"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "14.38"
$ws.Range("D14").Value = "21.43"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("E15").Value = "  +6.05%  "
$ws.Range("D16").Value = "0.764"
$ws.Range("D17").Value = "2.044.06"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "37.715.37"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "69.32"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "222.33"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +2.87%  "
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "18.75"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "1.29"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  +8.09%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "6.51"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("E37").Value = "  +4.19%  "
$ws.Range("E38").Value = "  +6.94%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "18.11"
$ws.Range("E40").Value = "  +7.90%  "
$ws.Range("D41").Value = "1.531.63"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "97.71"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "4.14"
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "6.99"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "2.233.86"
$ws.Range("E51").Value = "  +0.70%  "
